$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns D (Fecha, serial date), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
# for rows 2-12. These represent a row-wise re-shuffle of the weekly data.

$rows = @{
    2  = @{ D = 44432; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    3  = @{ D = 44421; J = 20; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí";  P = 600 }
    4  = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí";  P = 560 }
    5  = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    6  = @{ D = 44449; J = 30; K = 16000; L = 16000; M = 16000; O = "Provincia de Limarí";  P = 640 }
    7  = @{ D = 44446; J = 15; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí";  P = 520 }
    8  = @{ D = 44340; J = 25; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí";  P = 600 }
    9  = @{ D = 44453; J = 55; K = 14000; L = 15000; M = 14455; O = "Provincia de Limarí";  P = 578 }
    10 = @{ D = 44418; J = 12; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí";  P = 600 }
    11 = @{ D = 44376; J = 15; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí";  P = 480 }
    12 = @{ D = 44425; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí";  P = 560 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value = $vals.D
    $ws.Cells.Item($r, 10).Value = $vals.J
    $ws.Cells.Item($r, 11).Value = $vals.K
    $ws.Cells.Item($r, 12).Value = $vals.L
    $ws.Cells.Item($r, 13).Value = $vals.M
    $ws.Cells.Item($r, 15).Value = $vals.O
    $ws.Cells.Item($r, 16).Value = $vals.P
}
